# Apply the 2021 income fixes:
#  - correct RDLpc/RPLpc values (they were stored in thousands with a stray
#    decimal point, e.g. 11.178 instead of 11178) by multiplying by 1000
#  - the RPLpc (E) column additionally loses its inherited wrap-style during
#    the correction, ending up with the workbook's default "Normal" style
#  - add a small KPI-analysis marker cell (formatted, no value) at I4
#  - update the saved selection/scroll position
#  - set the print page to A4 portrait

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 75; $r++) {
  $c = $ws.Cells.Item($r, 3)
  $c.Value = $c.Value() * 1000

  $e = $ws.Cells.Item($r, 5)
  $e.Value = $e.Value() * 1000
}

# C75 had inherited the bold "totals row" style; the correction pass
# re-applies the plain wrapped style that the rest of column C uses.
$ws.Range("C75").Style = "Normal"
$ws.Range("C75").WrapText = $true
$ws.Range("C75").VerticalAlignment = -4108

# The whole RPLpc column (E) ends up unstyled (default "Normal" style).
$ws.Range("E2:E75").Style = "Normal"

# New KPI-analysis helper cell: formatted (underlined font, matching the
# sheet's footnote style) but left empty.
$ws.Range("I4").Font.Underline = 2

# Row heights grow slightly with the new default row height (14.5 -> 15pt).
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 45
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 45
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 30
$ws.Rows.Item(21).RowHeight = 45
$ws.Rows.Item(23).RowHeight = 45
$ws.Rows.Item(25).RowHeight = 30
$ws.Rows.Item(26).RowHeight = 45
$ws.Rows.Item(27).RowHeight = 45
$ws.Rows.Item(28).RowHeight = 30
$ws.Rows.Item(29).RowHeight = 45
$ws.Rows.Item(32).RowHeight = 30
$ws.Rows.Item(33).RowHeight = 60
$ws.Rows.Item(34).RowHeight = 30
$ws.Rows.Item(37).RowHeight = 30
$ws.Rows.Item(39).RowHeight = 30
$ws.Rows.Item(40).RowHeight = 45
$ws.Rows.Item(42).RowHeight = 30
$ws.Rows.Item(45).RowHeight = 45
$ws.Rows.Item(47).RowHeight = 30
$ws.Rows.Item(48).RowHeight = 30
$ws.Rows.Item(49).RowHeight = 30
$ws.Rows.Item(51).RowHeight = 30
$ws.Rows.Item(53).RowHeight = 30
$ws.Rows.Item(54).RowHeight = 30
$ws.Rows.Item(56).RowHeight = 30
$ws.Rows.Item(58).RowHeight = 30
$ws.Rows.Item(59).RowHeight = 30
$ws.Rows.Item(60).RowHeight = 30
$ws.Rows.Item(61).RowHeight = 30
$ws.Rows.Item(63).RowHeight = 30
$ws.Rows.Item(65).RowHeight = 45
$ws.Rows.Item(67).RowHeight = 45
$ws.Rows.Item(68).RowHeight = 60
$ws.Rows.Item(69).RowHeight = 30
$ws.Rows.Item(70).RowHeight = 75
$ws.Rows.Item(71).RowHeight = 30
$ws.Rows.Item(72).RowHeight = 45
$ws.Rows.Item(73).RowHeight = 45
$ws.Rows.Item(74).RowHeight = 30

# Update the view: scroll back to the top-left and move the active selection.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 1
$aw.ScrollColumn = 1
$ws.Range("G4").Select() | Out-Null

# Print setup: A4, portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
